$d = $word.ActiveDocument

# The paragraph reads (as separate runs): "...SCEC-VDO is a" + " " + "3D software...".
# Replace the "SCEC-VDO is a" run's text with "SCEC-VDO is an opensource" so the
# sentence becomes "...SCEC-VDO is an opensource 3D software...".
$d.Content.Find.Execute("SCEC-VDO is a", $true, $false, $false, $false, $false, $true, 1, $false, "SCEC-VDO is an opensource", 2)
